$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text format is preserved for Price (D) and Volume (E) columns
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.862.57'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '1.564.45'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '205.95'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '0.485'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '21.81'
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.786.49'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").Value = '1.565.83'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '3.73'
$ws.Range("E14").Value = '  -1.15%  '
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '26.863.65'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = '61.33'
$ws.Range("E17").Value = '  -2.41%  '
$ws.Range("D18").Value = '215.66'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '7.38'
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").Value = '2.00'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '153.50'
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("D27").Value = '14.97'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  -3.44%  '
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").Value = '1.402.12'
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").Value = '0.916'
$ws.Range("E37").Value = '  -3.68%  '
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").Value = '0.531'
$ws.Range("E39").Value = '  +3.10%  '
$ws.Range("D40").Value = '0.812'
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '0.994'
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("D43").Value = '5.57'
$ws.Range("E43").Value = '  +6.87%  '
$ws.Range("D44").Value = '1.80'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("D46").Value = '63.77'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = '1.700.08'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").Value = '86.61'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("D49").Value = '0.0503'
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").Value = '0.0₇0971'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '0.0952'
$ws.Range("E51").Value = '  +1.43%  '
